$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the newly-reported 2023 fiscal-year row of sectoral value-added data
# (Table 310-34101) right after the existing 2022 row (row 44).
$ws.Range("A45").Value = 2023
$ws.Range("B45").Value = 1131
$ws.Range("C45").Value = 28206
$ws.Range("D45").Value = 32794
$ws.Range("E45").Value = 126872
$ws.Range("F45").Value = 429774
$ws.Range("G45").Value = 15287
$ws.Range("H45").Value = 65619
$ws.Range("I45").Value = 16590
$ws.Range("J45").Value = 44372
$ws.Range("K45").Value = 44975
$ws.Range("L45").Value = 41985
$ws.Range("M45").Value = 71646
$ws.Range("N45").Value = 6821
$ws.Range("O45").Value = 624651
$ws.Range("P45").Value = 102495
$ws.Range("Q45").Value = 92352
$ws.Range("R45").Value = 9228
$ws.Range("S45").Value = 44685
$ws.Range("T45").Value = 56779
$ws.Range("U45").Value = 153715
$ws.Range("V45").Value = 320133
$ws.Range("W45").Value = 289981
$ws.Range("X45").Value = 296372

# Scroll/select back to the top of the sheet, as the author left it.
$ws.Range("A2").Select()
